$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the BOM table ("Tabelle1") to match the revised schematic ---

# Rename the last column header (cell F1) from "Unrat Best Nr." to
# "Conrad Best Nr." - writing the header cell also keeps the table's
# column name (xl/tables/table1.xml) in sync.
$ws.Range("F1").Value = "Conrad Best Nr."

# Row 3 (4,7u ceramic caps): one more capacitor (C51) was added at this
# position, so the count goes from 2 to 3 and the position list grows.
$ws.Range("A3").Value = 3
$ws.Range("D3").Value = "C1, C2, C51"

# Row 11 (100n caps row): two more capacitors (up to C50 instead of C48)
# were added, so the count goes from 42 to 44 and the range label updates.
$ws.Range("A11").Value = 44
$ws.Range("D11").Value = "C7, …, C50"

# Row 16 (220R resistor): value used to be the bare number 220; now it is
# stored as the text "220R".
$ws.Range("C16").Value = "220R"

# The old row 17 (R, 470, position R45, no order number) is obsolete and
# gets removed entirely - everything below it shifts up by one row.
$ws.Rows(17).Delete()

# Two new rows get appended to the table for the newly added components.
$lo = $ws.ListObjects.Item(1)
$lo.ListRows.Add() | Out-Null
$lo.ListRows.Add() | Out-Null

# New row 18: NPN transistors.
$ws.Range("A18").Value = 2
$ws.Range("B18").Value = "Q"
$ws.Range("C18").Value = "NPN Transistor"
$ws.Range("D18").Value = "Q1, Q2"
$ws.Range("E18").Value = "BCX 19 SMD"

# New row 19: additional 10k resistors.
$ws.Range("A19").Value = 2
$ws.Range("B19").Value = "R"
$ws.Range("C19").Value = "10k"
$ws.Range("D19").Value = "R45, R46, R47, R48"

# Match the author's final cursor position in the sheet.
$ws.Range("B3").Select()

Write-Output "BOM updated"
